$d = $word.ActiveDocument

# The final (empty) paragraph in the document currently just holds the
# "_GoBack" bookmark. Turn it into another bulleted "ListParagraph" item,
# continuing the same numbered/bulleted list (numId=1) that the preceding
# dialog lines use.

$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Item($paragraphs.Count)
$prevPara = $paragraphs.Item($paragraphs.Count - 1)

# Match the "List Paragraph" style used by the rest of the list items.
$lastPara.Style = "List Paragraph"

# Continue the bullet list already applied to the previous paragraph so the
# new paragraph reuses the same abstract numbering definition / numId
# instead of minting a brand-new list.
$listTemplate = $prevPara.Range.ListFormat.ListTemplate
$lastPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 0, $false, 1)
